$wb = $excel.ActiveWorkbook

# --- Table2: add row-21 averages (freq occurrence summary row) ---
$ws2 = $wb.Worksheets.Item("Table2")
$ws2.Activate()

$ws2.Range("B21").Formula = "=AVERAGE(B5:B20)"
$ws2.Range("C21:I21").Formula = "=AVERAGE(C5:C20)"

# Leave the active selection on E21, as in the edited workbook
$ws2.Range("E21").Select()

# --- Table3: move the active selection/cell ---
$ws3 = $wb.Worksheets.Item("Table3")
$ws3.Activate()
$ws3.Range("O18").Select()
